$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 190, shifting existing rows 190:263 down to 191:264.
$ws.Rows.Item(190).Insert("xlShiftDown")

# Populate the newly inserted row 190 with the new observation.
# Columns A,B,C,E,F,G,H,I,J,N,O,Q,R repeat the values of the (now shifted) row 191,
# only D (fecha), K/L/M (precios) and P (precio $/Kg) differ.
$ws.Cells.Item(190, 1).Value  = 7
$ws.Cells.Item(190, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(190, 3).Value  = "Ñuble"
$ws.Cells.Item(190, 4).Value  = 44900
$ws.Cells.Item(190, 5).Value  = 16
$ws.Cells.Item(190, 6).Value  = 100112017
$ws.Cells.Item(190, 7).Value  = "Apio"
$ws.Cells.Item(190, 8).Value  = "Americana (o)"
$ws.Cells.Item(190, 9).Value  = "Primera"
$ws.Cells.Item(190, 10).Value = 120
$ws.Cells.Item(190, 11).Value = 9000
$ws.Cells.Item(190, 12).Value = 9500
$ws.Cells.Item(190, 13).Value = 9250
$ws.Cells.Item(190, 14).Value = "$/docena de matas"
$ws.Cells.Item(190, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(190, 16).Value = 1542
$ws.Cells.Item(190, 17).Value = 6
$ws.Cells.Item(190, 18).Value = "Hortaliza"
